$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the old row 805 (the "Spiders and Knives (1)" entry)
# to hold the two new gender-restriction TLK string entries.
$ws.Rows("805:806").Insert()

# The freshly inserted rows inherit the fill/format of the row above (804); the
# real edit used the plain/default string-row style instead, so copy that
# formatting in from a representative unstyled row (row 2) before filling values.
$ws.Range("A2:B2").Copy()
$ws.Range("A805:B806").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A805").Value = 6610803
$ws.Range("B805").Value = "Required: Female"

$ws.Range("A806").Value = 6610804
$ws.Range("B806").Value = "Required: Male"

# New comment describing the purpose of these two rows.
$ws.Range("A805").AddComment("Item restrictions")

# Reflect the author's final on-screen selection / scroll position.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 775
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B807").Select()
